$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New expense rows (41-44) -------------------------------------------
# Row 41: Legacy (Remeras Trabajo) - 1935 in October (column D)
$ws.Range("A41").Value = "Legacy (Remeras Trabajo)"
$ws.Range("D41").Value = 1935

# Row 42: Filtro de Agua - 329 every month, October (D) through a new
# column O (next month added to the tracker)
$ws.Range("A42").Value = "Filtro de Agua"
$ws.Range("N42").Copy()
$ws.Range("O42").PasteSpecial(-4122)
$ws.Range("D42:O42").Value = 329

# Row 43: Legacy (Jean para salir) - 1345 in October and November (D:E)
$ws.Range("A43").Value = "Legacy (Jean para salir)"
$ws.Range("D43").Value = 1345
$ws.Range("E43").Value = 1345

# Row 44: Recarga celular - 100 in October (D)
$ws.Range("A44").Value = "Recarga celular"
$ws.Range("D44").Value = 100

# --- Gastos Comunes correction ------------------------------------------
# OCA bill for October revised down from 2500 to 2200
$ws.Range("D55").Value = 2200

# --- Refresh the totals so Excel records them as shared formulas --------
# (mirrors re-filling the formula across the row, which is how these
# shared-formula blocks originally got created)
$ws.Range("B50:N50").Formula = "=SUM(B8:B49)"
$ws.Range("B56:N56").Formula = "=SUM(B52:B55)"
$ws.Range("B61:N61").Formula = "=SUM(B50+B56+B58)"

# --- View state: selection now parked on the new bottom-line total ------
$ws.Range("D61").Select()
